$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-138). The automatic update bumps that date from 45203 to 45205
# for every row.
$ws.Range("C2:C138").Value = 45205
